# Update the build timestamp embedded in the workbook's version strings
# from "February 03 2026 17.29.55 EST" to "February 03 2026 18.05.36 EST".

$wb = $excel.ActiveWorkbook

$oldStamp = "February 03 2026 17.29.55 EST"
$newStamp = "February 03 2026 18.05.36 EST"

# --- "About" sheet: A2 (Version line) and A6 (Recommended Citation line) ---
$wsAbout = $wb.Worksheets.Item("About")

$a2 = $wsAbout.Range("A2").Value()
$wsAbout.Range("A2").Value = $a2.Replace($oldStamp, $newStamp)

$a6 = $wsAbout.Range("A6").Value()
$wsAbout.Range("A6").Value = $a6.Replace($oldStamp, $newStamp)

# --- "Boundaries and methane sources" sheet: S2:S25 (build_version column) ---
$wsData = $wb.Worksheets.Item("Boundaries and methane sources")

for ($row = 2; $row -le 25; $row++) {
    $cell = $wsData.Cells.Item($row, 19)  # column S = 19
    $val = $cell.Value()
    if ($val -ne $null) {
        $cell.Value = $val.Replace($oldStamp, $newStamp)
    }
}
